# CRMS-2130 add one more column
# Insert a new "Partner Type" column right after "Partner Status" (column AD).
# This shifts the old columns AE:AK one to the right (AF:AL) and adds the new
# header text in row 1 plus the corresponding merge-field placeholder in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at AE, shifting AE:AK -> AF:AL
$ws.Columns("AE:AE").Insert()

# Populate the new column's header (row 1) and placeholder (row 2)
$ws.Range("AE1").Value = "Partner Type"
$ws.Range("AE2").Value = "{excel_data_line_item:partner_type}"

# Update the active selection to match the edited cell
$ws.Range("AE2").Select()
